# remove outliers should be fixed
# - Update the "data" sheet's outlier value (C5) from 1 to 30
# - Make "data" the active sheet/tab, with selection on G5
# - "setup" sheet is no longer the selected/active tab

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")

# Fix the outlier value
$wsData.Range("C5").Value = 30

# Activate the data sheet and set its selection
$wsData.Activate()
$wsData.Range("G5").Select()
